# lab4/results.xlsx — "results of ex2 30 tasks"
#
# The workbook originally held one sheet ("Feuil1") with results for the
# first experiment. This edit:
#   1. Renames "Feuil1" -> "experiment1"
#   2. Duplicates it into a second sheet "experiment2" (placed after the
#      first), which becomes the active sheet/tab.
#   3. On BOTH sheets, the "var" header text is replaced by "pold" (column
#      B) / "moyenne" (columns F and K) — "var" ends up unused and is
#      dropped from the shared-string table on save.
#   4. On "experiment2", the raw per-run cost/time samples (columns C/D)
#      are cleared out (no data collected for this run -> #DIV/0! on the
#      AVERAGE formulas), while the raw pold/cost samples (columns H/I) are
#      replaced with the new 30-task experiment measurements.
#   5. Selections are adjusted to match where each sheet's user was last
#      positioned (K2 on experiment1, H26 on experiment2).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1/2: duplicate the sheet, placing the copy right after the original.
# Worksheets.Add() would insert *before*, so we use Copy(After:=) instead;
# Excel makes the new copy the active sheet automatically (-> activeTab=1).
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)

$ws1.Name = "experiment1"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "experiment2"

# --- 3: fix up the header row text on both sheets ("var" -> retired).
foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("B1").Value = "pold"
    $ws.Range("F1").Value = "moyenne"
    $ws.Range("K1").Value = "moyenne"
}

# --- 4: experiment2 holds the new 30-task run.
# No cost/time samples were recorded -> clear C2:D25 (AVERAGE formulas in
# E/F will naturally fall back to #DIV/0!).
$ws2.Range("C2:D25").ClearContents()

# New pold/cost raw samples for the 30-task experiment (columns H/I).
$samples = @(
    [PSCustomObject]@{ Row = 2;  H = 18036; I = 70139  },
    [PSCustomObject]@{ Row = 3;  H = 23408; I = 33315  },
    [PSCustomObject]@{ Row = 4;  H = 21877; I = 22549  },
    [PSCustomObject]@{ Row = 5;  H = 18718; I = 99876  },
    [PSCustomObject]@{ Row = 6;  H = 18543; I = 95819  },
    [PSCustomObject]@{ Row = 7;  H = 18802; I = 53744  },
    [PSCustomObject]@{ Row = 8;  H = 13486; I = 159510 },
    [PSCustomObject]@{ Row = 9;  H = 18314; I = 66751  },
    [PSCustomObject]@{ Row = 10; H = 17820; I = 65164  },
    [PSCustomObject]@{ Row = 11; H = 16835; I = 78205  },
    [PSCustomObject]@{ Row = 12; H = 19209; I = 85517  },
    [PSCustomObject]@{ Row = 13; H = 18758; I = 69015  },
    [PSCustomObject]@{ Row = 14; H = 18257; I = 61657  },
    [PSCustomObject]@{ Row = 15; H = 16379; I = 110294 },
    [PSCustomObject]@{ Row = 16; H = 14885; I = 302021 },
    [PSCustomObject]@{ Row = 17; H = 12543; I = 302020 },
    [PSCustomObject]@{ Row = 18; H = 16485; I = 131320 },
    [PSCustomObject]@{ Row = 19; H = 19882; I = 187975 },
    [PSCustomObject]@{ Row = 20; H = 20159; I = 54902  },
    [PSCustomObject]@{ Row = 21; H = 19929; I = 62572  },
    [PSCustomObject]@{ Row = 22; H = 21930; I = 16379  },
    [PSCustomObject]@{ Row = 23; H = 15286; I = 167206 },
    [PSCustomObject]@{ Row = 24; H = 18907; I = 81267  },
    [PSCustomObject]@{ Row = 25; H = 13004; I = 150093 }
)

foreach ($s in $samples) {
    $ws2.Cells.Item($s.Row, 8).Value = $s.H
    $ws2.Cells.Item($s.Row, 9).Value = $s.I
}

# --- 5: restore each sheet's last-used selection.
$ws1.Range("K2").Select()
$ws2.Range("H26").Select()
